# Update the "want to go" counts (column F) on both the "展览" sheet and its
# "全部类型" duplicate. Sheet "展览" gets F5 = 8496, while "全部类型" gets F5 = 8497;
# every other updated row is identical between the two sheets.

$wb = $excel.ActiveWorkbook

$commonUpdates = @{
    3  = 506
    4  = 438
    7  = 1503
    8  = 175
    11 = 245
    12 = 384
    16 = 127
    19 = 1224
    20 = 175
    23 = 90
    24 = 121
    25 = 67
    27 = 100
}

$sheetSpecificF5 = @{
    "展览"     = 8496
    "全部类型" = 8497
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $commonUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $commonUpdates[$row]
    }

    $ws.Cells.Item(5, 6).Value = $sheetSpecificF5[$sheetName]
}
